{"js": "// Lattice-multiplication worksheet refresh: every exercise cell in the\n// (single) table gets a new \"A x B\" problem. The rest of each cell\n// (the spaced-out digits of B, the dashed rule, and the two boxed\n// digits of A) is fully derived from \"A x B\", so we recompute all five\n// lines for every cell instead of trying to patch individual numbers.\n\n// New \"A x B\" values, in row-major order (row 0 = top row, col 0 = left column).\nconst NEW_PROBLEMS = [\n  [\"55 x 43\", \"38 x 84\", \"60 x 16\"],\n  [\"56 x 98\", \"35 x 44\", \"18 x 15\"],\n  [\"53 x 97\", \"69 x 72\", \"37 x 73\"],\n  [\"77 x 71\", \"88 x 20\", \"41 x 57\"],\n  [\"27 x 81\", \"49 x 55\", \"87 x 38\"],\n];\n\n// Build the 5-line lattice text (joined with vertical-tab \\u000b, the\n// character Word uses in-memory for a line break / <w:br/>) for a given\n// \"A x B\" problem string.\nfunction buildCellText(problem) {\n  const [aStr, , bStr] = problem.split(\" \");\n  const a = aStr.trim();\n  const b = bStr.trim();\n  const line1 = problem;\n  const line2 = `  ${b[0]}    ${b[1]}`;\n  const line3 = \"  ----\";\n  const line4 = `${a[0]}|    |`;\n  const line5 = `${a[1]}|    |`;\n  return [line1, line2, line3, line4, line5].join(\"\\u000b\");\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length && r < NEW_PROBLEMS.length; r++) {\n  const cells = rows.items[r].cells.items;\n  const rowValues = NEW_PROBLEMS[r];\n  for (let c = 0; c < cells.length && c < rowValues.length; c++) {\n    const cell = cells[c];\n    const paragraph = cell.body.paragraphs.getFirst();\n    const range = paragraph.getRange(\"Whole\");\n    range.insertText(buildCellText(rowValues[c]), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet refresh: every exercise cell in the\n# (single) table gets a new \"A x B\" problem. The rest of each cell (the\n# spaced-out digits of B, the dashed rule, and the two boxed digits of\n# A) is fully derived from \"A x B\", so we recompute all five lines for\n# every cell instead of trying to patch individual numbers.\n\n$d = $word.ActiveDocument\n\n# New \"A x B\" values, in row-major order (row 0 = top row, col 0 = left column).\n$newProblems = @(\n    @(\"55 x 43\", \"38 x 84\", \"60 x 16\"),\n    @(\"56 x 98\", \"35 x 44\", \"18 x 15\"),\n    @(\"53 x 97\", \"69 x 72\", \"37 x 73\"),\n    @(\"77 x 71\", \"88 x 20\", \"41 x 57\"),\n    @(\"27 x 81\", \"49 x 55\", \"87 x 38\")\n)\n\n$vtab = [char]11\n\nfunction Build-CellText($problem) {\n    $parts = $problem -split \" x \"\n    $a = $parts[0]\n    $b = $parts[1]\n    $a0 = $a.Substring(0,1)\n    $a1 = $a.Substring(1,1)\n    $b0 = $b.Substring(0,1)\n    $b1 = $b.Substring(1,1)\n    $line1 = $problem\n    $line2 = \"  {0}    {1}\" -f $b0, $b1\n    $line3 = \"  ----\"\n    $line4 = \"{0}|    |\" -f $a0\n    $line5 = \"{0}|    |\" -f $a1\n    return ($line1 + $vtab + $line2 + $vtab + $line3 + $vtab + $line4 + $vtab + $line5)\n}\n\n$table = $d.Tables.Item(1)\n$rowCount = $newProblems.Count\nfor ($r = 0; $r -lt $rowCount; $r++) {\n    $rowValues = $newProblems[$r]\n    for ($c = 0; $c -lt $rowValues.Count; $c++) {\n        $cell = $table.Cell($r + 1, $c + 1)\n        $range = $cell.Range\n        # Trim the trailing cell-mark/paragraph-mark pair off the end of\n        # the cell range so we only overwrite the visible text.\n        $range.End = $range.End - 2\n        $range.Text = Build-CellText $rowValues[$c]\n    }\n}\n\nWrite-Output \"done\"\n"}
